# "removed launch hotkey + modified data"
#
# Row 6 / column A used to hold the shared string
# "Nous mobilisons une équipe de Déménageurs Qualifiés" (which also gets
# dropped from the shared-strings table entirely). It is replaced with a
# fresh ad-posting URL, which becomes a brand-new shared string appended
# at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A6").Value = "leboncoin.fr/deposer-une-annonce"

# Selection moved from D8 to A7, and the window was resized/repositioned
# (cosmetic view state saved along with the workbook).
$ws.Range("A7").Select() | Out-Null

$win = $wb.Windows.Item(1)
$win.WindowState = -4143  # xlNormal
$win.Left = -120
$win.Top = -120
$win.Width = 24240
$win.Height = 13140 | Out-Null
